$d = $word.ActiveDocument

# The first paragraph currently reads "This is a Microsoft word document."
# We need to append the text " (Changed main)" as three additional runs:
#   " (", "Changed main", ")"
$para = $d.Paragraphs(1).Range
$end = $para.End - 1   # before the paragraph mark

$r1 = $d.Range($end, $end)
$r1.InsertAfter(" (")

$end = $end + 2
$r2 = $d.Range($end, $end)
$r2.InsertAfter("Changed main")

$end = $end + 12
$r3 = $d.Range($end, $end)
$r3.InsertAfter(")")
